# chore: update Sheets via scheduled runner
# Refreshes the market-board derived pricing/profit figures (columns H:N)
# on the per-crafting-job leve profit sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1322.7
$ws.Cells.Item(15, 9).Value = 1322.7
$ws.Cells.Item(15, 11).Value = 3968.1
$ws.Cells.Item(15, 13).Value = -3799.1
$ws.Cells.Item(112, 8).Value = 2648.5264
$ws.Cells.Item(112, 9).Value = 1777
$ws.Cells.Item(112, 10).Value = 2696.9443
$ws.Cells.Item(112, 11).Value = 5331
$ws.Cells.Item(112, 12).Value = 8090.8329
$ws.Cells.Item(112, 13).Value = -4223
$ws.Cells.Item(112, 14).Value = -10306.8329
$ws.Cells.Item(137, 8).Value = 4537.9707
$ws.Cells.Item(137, 9).Value = 6074.5
$ws.Cells.Item(137, 10).Value = 4065.1924
$ws.Cells.Item(137, 11).Value = 18223.5
$ws.Cells.Item(137, 12).Value = 12195.5772
$ws.Cells.Item(137, 13).Value = -15673.5
$ws.Cells.Item(137, 14).Value = -17295.5772
$ws.Cells.Item(138, 8).Value = 8253.485000000001
$ws.Cells.Item(138, 9).Value = 3739.4167
$ws.Cells.Item(138, 10).Value = 10608.652
$ws.Cells.Item(138, 11).Value = 11218.2501
$ws.Cells.Item(138, 12).Value = 31825.956
$ws.Cells.Item(138, 13).Value = -6078.250100000001
$ws.Cells.Item(138, 14).Value = -42105.956
$ws.Cells.Item(141, 8).Value = 2872.4666
$ws.Cells.Item(141, 9).Value = 2863.3572
$ws.Cells.Item(141, 11).Value = 8590.071599999999
$ws.Cells.Item(141, 13).Value = -3410.071599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 11555.167
$ws.Cells.Item(45, 9).Value = 3999
$ws.Cells.Item(45, 10).Value = 15333.25
$ws.Cells.Item(45, 11).Value = 3999
$ws.Cells.Item(45, 12).Value = 15333.25
$ws.Cells.Item(45, 13).Value = -3622
$ws.Cells.Item(45, 14).Value = -16087.25
$ws.Cells.Item(61, 8).Value = 4836.409
$ws.Cells.Item(61, 9).Value = 2681.4167
$ws.Cells.Item(61, 10).Value = 7422.4
$ws.Cells.Item(61, 11).Value = 2681.4167
$ws.Cells.Item(61, 12).Value = 7422.4
$ws.Cells.Item(61, 13).Value = -2469.4167
$ws.Cells.Item(61, 14).Value = -7846.4
$ws.Cells.Item(74, 8).Value = 316316.66
$ws.Cells.Item(74, 9).Value = 437650.47
$ws.Cells.Item(74, 10).Value = 6241.3335
$ws.Cells.Item(74, 11).Value = 437650.47
$ws.Cells.Item(74, 12).Value = 6241.3335
$ws.Cells.Item(74, 13).Value = -436776.47
$ws.Cells.Item(74, 14).Value = -7989.3335
$ws.Cells.Item(77, 8).Value = 316316.66
$ws.Cells.Item(77, 9).Value = 437650.47
$ws.Cells.Item(77, 10).Value = 6241.3335
$ws.Cells.Item(77, 11).Value = 2188252.35
$ws.Cells.Item(77, 12).Value = 31206.6675
$ws.Cells.Item(77, 13).Value = -2183884.35
$ws.Cells.Item(77, 14).Value = -39942.6675
$ws.Cells.Item(132, 8).Value = 6991.564
$ws.Cells.Item(132, 9).Value = 2469.1304
$ws.Cells.Item(132, 10).Value = 13492.5625
$ws.Cells.Item(132, 11).Value = 7407.3912
$ws.Cells.Item(132, 12).Value = 40477.6875
$ws.Cells.Item(132, 13).Value = -4877.3912
$ws.Cells.Item(132, 14).Value = -45537.6875
$ws.Cells.Item(136, 8).Value = 4836.409
$ws.Cells.Item(136, 9).Value = 2681.4167
$ws.Cells.Item(136, 10).Value = 7422.4
$ws.Cells.Item(136, 11).Value = 8044.250100000001
$ws.Cells.Item(136, 12).Value = 22267.2
$ws.Cells.Item(136, 13).Value = -5494.250100000001
$ws.Cells.Item(136, 14).Value = -27367.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 21672.34
$ws.Cells.Item(134, 9).Value = 2006.8334
$ws.Cells.Item(134, 10).Value = 96758.82000000001
$ws.Cells.Item(134, 11).Value = 6020.5002
$ws.Cells.Item(134, 12).Value = 290276.46
$ws.Cells.Item(134, 13).Value = -3485.5002
$ws.Cells.Item(134, 14).Value = -295346.46

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 79397.5
$ws.Cells.Item(20, 10).Value = 79397.5
$ws.Cells.Item(20, 12).Value = 79397.5
$ws.Cells.Item(20, 14).Value = -79869.5
$ws.Cells.Item(30, 8).Value = 79397.5
$ws.Cells.Item(30, 10).Value = 79397.5
$ws.Cells.Item(30, 12).Value = 79397.5
$ws.Cells.Item(30, 14).Value = -79579.5
$ws.Cells.Item(31, 8).Value = 5324.737
$ws.Cells.Item(31, 9).Value = 2627
$ws.Cells.Item(31, 10).Value = 5933.9033
$ws.Cells.Item(31, 11).Value = 2627
$ws.Cells.Item(31, 12).Value = 5933.9033
$ws.Cells.Item(31, 13).Value = -2332
$ws.Cells.Item(31, 14).Value = -6523.9033
$ws.Cells.Item(34, 8).Value = 5324.737
$ws.Cells.Item(34, 9).Value = 2627
$ws.Cells.Item(34, 10).Value = 5933.9033
$ws.Cells.Item(34, 11).Value = 2627
$ws.Cells.Item(34, 12).Value = 5933.9033
$ws.Cells.Item(34, 13).Value = -2425
$ws.Cells.Item(34, 14).Value = -6337.9033
$ws.Cells.Item(99, 8).Value = 4041.6206
$ws.Cells.Item(99, 9).Value = 3052.647
$ws.Cells.Item(99, 10).Value = 5442.6665
$ws.Cells.Item(99, 11).Value = 3052.647
$ws.Cells.Item(99, 12).Value = 5442.6665
$ws.Cells.Item(99, 13).Value = -1554.647
$ws.Cells.Item(99, 14).Value = -8438.666499999999
$ws.Cells.Item(126, 8).Value = 4041.6206
$ws.Cells.Item(126, 9).Value = 3052.647
$ws.Cells.Item(126, 10).Value = 5442.6665
$ws.Cells.Item(126, 11).Value = 9157.940999999999
$ws.Cells.Item(126, 12).Value = 16327.9995
$ws.Cells.Item(126, 13).Value = -6687.940999999999
$ws.Cells.Item(126, 14).Value = -21267.9995
$ws.Cells.Item(128, 8).Value = 79397.5
$ws.Cells.Item(128, 10).Value = 79397.5
$ws.Cells.Item(128, 12).Value = 79397.5
$ws.Cells.Item(128, 14).Value = -89357.5
$ws.Cells.Item(134, 8).Value = 560188.9399999999
$ws.Cells.Item(134, 9).Value = 4127.4546
$ws.Cells.Item(134, 10).Value = 1433999.9
$ws.Cells.Item(134, 11).Value = 12382.3638
$ws.Cells.Item(134, 12).Value = 4301999.699999999
$ws.Cells.Item(134, 13).Value = -9847.363799999999
$ws.Cells.Item(134, 14).Value = -4307069.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(42, 8).Value = 2334.6667
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(68, 8).Value = 3315.8948
$ws.Cells.Item(68, 9).Value = 2460.4167
$ws.Cells.Item(68, 10).Value = 4782.4287
$ws.Cells.Item(68, 11).Value = 7381.250100000001
$ws.Cells.Item(68, 12).Value = 14347.2861
$ws.Cells.Item(68, 13).Value = -6570.250100000001
$ws.Cells.Item(68, 14).Value = -15969.2861
$ws.Cells.Item(69, 8).Value = 9528.571
$ws.Cells.Item(69, 9).Value = 6700
$ws.Cells.Item(69, 11).Value = 20100
$ws.Cells.Item(69, 13).Value = -19289
$ws.Cells.Item(71, 8).Value = 3315.8948
$ws.Cells.Item(71, 9).Value = 2460.4167
$ws.Cells.Item(71, 10).Value = 4782.4287
$ws.Cells.Item(71, 11).Value = 22143.7503
$ws.Cells.Item(71, 12).Value = 43041.85830000001
$ws.Cells.Item(71, 13).Value = -18087.7503
$ws.Cells.Item(71, 14).Value = -51153.85830000001
$ws.Cells.Item(72, 8).Value = 9528.571
$ws.Cells.Item(72, 9).Value = 6700
$ws.Cells.Item(72, 11).Value = 60300
$ws.Cells.Item(72, 13).Value = -56244
$ws.Cells.Item(113, 8).Value = 2179663
$ws.Cells.Item(113, 9).Value = 6174089
$ws.Cells.Item(113, 10).Value = 885.1818
$ws.Cells.Item(113, 11).Value = 18522267
$ws.Cells.Item(113, 12).Value = 2655.5454
$ws.Cells.Item(113, 13).Value = -18520097
$ws.Cells.Item(113, 14).Value = -6995.5454
$ws.Cells.Item(42, 13).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1410.2
$ws.Cells.Item(22, 9).Value = 1783.6666
$ws.Cells.Item(22, 10).Value = 850
$ws.Cells.Item(22, 11).Value = 1783.6666
$ws.Cells.Item(22, 12).Value = 850
$ws.Cells.Item(22, 13).Value = -1488.6666
$ws.Cells.Item(22, 14).Value = -1440
$ws.Cells.Item(27, 8).Value = 1410.2
$ws.Cells.Item(27, 9).Value = 1783.6666
$ws.Cells.Item(27, 10).Value = 850
$ws.Cells.Item(27, 11).Value = 1783.6666
$ws.Cells.Item(27, 12).Value = 850
$ws.Cells.Item(27, 13).Value = -1676.6666
$ws.Cells.Item(27, 14).Value = -1064
$ws.Cells.Item(43, 8).Value = 5749699
$ws.Cells.Item(43, 9).Value = 5749699
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 5749699
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = -5749506
$ws.Cells.Item(43, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 357313.88
$ws.Cells.Item(136, 9).Value = 419427.88
$ws.Cells.Item(136, 10).Value = 208240.3
$ws.Cells.Item(136, 11).Value = 1258283.64
$ws.Cells.Item(136, 12).Value = 624720.8999999999
$ws.Cells.Item(136, 13).Value = -1255733.64
$ws.Cells.Item(136, 14).Value = -629820.8999999999
